$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New "Usable" column (H) content
#    H1 = header "Usable"
#    H2 = "No"   (subnet 0 - network address, not usable)
#    H3:H16 = "Yes"
#    H17 = "No"  (subnet 15 - broadcast address, not usable)
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Usable"
$ws.Range("H2").Value = "No"
$ws.Range("H3:H16").Value = "Yes"
$ws.Range("H17").Value = "No"

# Match the text/right-aligned formatting already used by the rest of the
# data rows (some of these H cells - rows 10-17 - did not exist before, so
# they'd otherwise pick up the plain default style instead).
$dataCol = $ws.Range("H2:H17")
$dataCol.NumberFormat = "@"
$dataCol.HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 2. Borders: thin grid around the whole table, thick rule under the header
#    row, and no rule directly above row 2 (it sits right under the thick
#    header rule so a plain thin rule there would double up).
# ---------------------------------------------------------------------------
$full = $ws.Range("A1:H17")
$full.Borders.LineStyle = 1
$full.Borders.Weight = 2

$row2 = $ws.Range("A2:H2")
$row2.Borders.Item(8).LineStyle = 0

$header = $ws.Range("A1:H1")
$header.Borders.Item(9).LineStyle = 1
$header.Borders.Item(9).Weight = 4

# ---------------------------------------------------------------------------
# 3. View state: zoom + active selection
# ---------------------------------------------------------------------------
$ws.Range("J17").Select()
$excel.ActiveWindow.Zoom = 130
